# Case and Fatality Demographics Data Updated
# Updates the "Fatalities by ..." sheets with refreshed counts and
# recalculated percentages, matching the 2021-10-29 data refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")

$wsAge.Range("B3").Value = 20
$wsAge.Range("B4").Value = 76
$wsAge.Range("B5").Value = 608
$wsAge.Range("B6").Value = 1946
$wsAge.Range("B7").Value = 4846
$wsAge.Range("B8").Value = 9231
$wsAge.Range("B9").Value = 7077
$wsAge.Range("B10").Value = 8391
$wsAge.Range("B11").Value = 8962
$wsAge.Range("B12").Value = 8506
$wsAge.Range("B13").Value = 20170
$wsAge.Range("B15").Value = 69846

$wsAge.Activate() | Out-Null
$wsAge.Range("B2:B15").Select() | Out-Null

# ---------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")

$wsGender.Range("B2").Value = 29201
$wsGender.Range("B3").Value = 40644
# B4 (Unknown) is unchanged
# Total (B5) becomes a plain value instead of =SUM(B2:B4)
$wsGender.Range("B5").Value = 69846

$wsGender.Activate() | Out-Null
$wsGender.Range("B2:B5").Select() | Out-Null

# ---------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

$wsRace.Range("B2").Value = 1296
$wsRace.Range("B3").Value = 7310
$wsRace.Range("B4").Value = 30565
$wsRace.Range("B5").Value = 412
$wsRace.Range("B6").Value = 30221
$wsRace.Range("B7").Value = 42
# B8 keeps its =SUM(B2:B7) formula and recalculates automatically

$wsRace.Activate() | Out-Null
$wsRace.Range("E16").Select() | Out-Null
